# Update the dSF (column F) values following a repull of the underlying
# data / push of all data and recalculation of the mean.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -1
$ws.Range("F4").Value  = -3
$ws.Range("F6").Value  = -2
$ws.Range("F8").Value  = -4
$ws.Range("F9").Value  = -5
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 5
$ws.Range("F13").Value = -5
